# Flight_Mission_Cycle.xlsx - "Fixed errors, added error warnings"
#
# The "Force_End" row label on the Writing sheet was renamed to "Force",
# and the sheet selection moved from J10 to I13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Writing")

# Rename the "Force_End" label (row 2, column A) to "Force".
$ws.Range("A2").Value = "Force"

# Make sure "Writing" is the active sheet, then update the saved selection
# to I13 (previously J10).
$ws.Activate()
$ws.Range("I13").Select()
